# Updated cryptos list on Wed Jun 26 17:08:16 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns for each coin row,
# and fixes the Stacks/Mantle row ordering (rows 42-43).
#
# Some refreshed Price values look like plain numbers (e.g. "0.468",
# "1.70"). The sheet stores prices as literal text (to preserve exact
# formatting such as trailing zeros / thousand separators), so for those
# cells we force the cell's NumberFormat to Text ("@") before assigning
# the value - otherwise Excel's COM layer would silently coerce the
# string into a floating point number and mangle the display text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.127.18"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "3.344.27"
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.32"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.79"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "3.346.03"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.468"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.44"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("E12").Value = "  -2.02%  "
$ws.Range("D13").Value = "3.925.58"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  -3.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.76"
$ws.Range("E16").Value = "  +1.65%  "
$ws.Range("D17").Value = "3.354.66"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "61.268.87"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.91"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.81"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.24"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.26"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("E23").Value = "  -3.91%  "
$ws.Range("D24").Value = "3.498.72"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.92"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000123"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.77"
$ws.Range("E28").Value = "  +9.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.46"
$ws.Range("E29").Value = "  -3.37%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.13"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.13"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("E34").Value = "  +0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.41"
$ws.Range("E35").Value = "  -0.39%  "
$ws.Range("E36").Value = "  -5.69%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.75"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "165.05"
$ws.Range("E38").Value = "  +2.16%  "
$ws.Range("E39").Value = "  -2.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0756"
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.766"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.70"
$ws.Range("E43").Value = "  -1.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.37"
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.36"
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.85"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.78"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.85"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").Value = "2.351.29"
$ws.Range("E50").Value = "  -0.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0259"
$ws.Range("E51").Value = "  -2.27%  "
